$d = $word.ActiveDocument

$replacements = @(
    @{old = "999÷5="; new = "733÷4="},
    @{old = "672÷2="; new = "160÷8="},
    @{old = "433÷2="; new = "214÷5="},
    @{old = "670÷7="; new = "408÷9="},
    @{old = "229÷7="; new = "355÷9="},
    @{old = "655÷3="; new = "655÷9="},
    @{old = "730÷4="; new = "928÷2="},
    @{old = "176÷5="; new = "987÷7="},
    @{old = "690÷5="; new = "288÷8="},
    @{old = "834÷7="; new = "250÷6="},
    @{old = "539÷2="; new = "607÷8="},
    @{old = "956÷6="; new = "566÷5="},
    @{old = "940÷9="; new = "243÷3="},
    @{old = "480÷8="; new = "232÷5="},
    @{old = "738÷4="; new = "457÷9="},
    @{old = "591÷7="; new = "791÷6="},
    @{old = "673÷6="; new = "316÷8="},
    @{old = "229÷2="; new = "948÷3="},
    @{old = "741÷9="; new = "804÷6="},
    @{old = "301÷8="; new = "825÷8="},
    @{old = "911÷7="; new = "881÷3="},
    @{old = "265÷2="; new = "910÷7="},
    @{old = "297÷2="; new = "726÷5="},
    @{old = "533÷2="; new = "522÷3="},
    @{old = "674÷9="; new = "654÷7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
